# Applies the commit "Tentando terminar API @TheoMesquita" edits to the
# "Modelo - Documento Basico - SPTECH" document.
#
# Summary of changes:
#  1. Merge the author name + registration number into a single run.
#  2. Rewrite the "Contexto do Negocio" paragraph text.
#  3. Remove one of the two blank paragraphs before "Objetivo".
#  4. Rewrite the "Objetivo" paragraph, split across several runs.
#  5. Rewrite the "Justificativa" paragraph, split across several runs.
#  6. Leave "Escopo" section untouched.
#  7. Rewrite the "Premissas e Restricoes" paragraph, split across several runs.

$d = $word.ActiveDocument

function Set-MultiRunParagraph {
    param(
        $doc,
        [string]$FindText,
        [string[]]$Segments
    )

    $r = $doc.Content
    $found = $r.Find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $FindText"
        return
    }

    $start = $r.Start
    $end = $r.End
    $target = $doc.Range($start, $end)
    $joined = [string]::Join("", $Segments)
    $target.Text = $joined

    # Force a run boundary between consecutive segments by toggling Bold
    # on/off over the just-inserted sub-range (identical formatting before
    # and after means the visible result is unaffected, but it keeps the
    # text split into discrete <w:r> elements instead of being coalesced
    # into a single run).
    $pos = $start
    for ($i = 0; $i -lt $Segments.Length; $i++) {
        $seg = $Segments[$i]
        $segStart = $pos
        $segEnd = $pos + $seg.Length
        if ($i -gt 0) {
            $sr = $doc.Range($segStart, $segEnd)
            $sr.Bold = 1
            $sr.Bold = 0
        }
        $pos = $segEnd
    }
}

# 1) "Theofilo Fernandes de Mesquita" + " - 01231042" -> single run
$d.Content.Find.Execute(
    "Theofilo Fernandes de Mesquita – 01231042", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Theofilo Fernandes de Mesquita – 01231042", 2) | Out-Null

# 2) "Contexto do Negocio" paragraph text rewrite
$d.Content.Find.Execute(
    "Site criado com o intuito de trazer ao público entretenimento como para assistir filmes, séries e animes!",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Esse site criado contendo conteudos para pessoas de todas as idades ele tem o proposito de mostrar como é um site de filmes como aqueles famosos como a Netflix, para posteriormente ser usado como um site real de filmes e series funcional.",
    2) | Out-Null

# 3) Remove one of the two blank paragraphs right before "Objetivo", and
#    one of the two blank paragraphs right before "Justificativa".
function Remove-BlankParagraphBeforeHeading {
    param($doc, [string]$HeadingText)

    $idx = 0
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text.TrimEnd([char]13) -eq $HeadingText) {
            $idx = $i
            break
        }
    }
    if ($idx -gt 1) {
        $doc.Paragraphs($idx - 1).Range.Delete() | Out-Null
    }
}

Remove-BlankParagraphBeforeHeading $d "Objetivo"
Remove-BlankParagraphBeforeHeading $d "Justificativa"

# 4) "Objetivo" paragraph content rewrite (split across 5 runs)
Set-MultiRunParagraph $d "Construir um site funcional de filmes, séries e animes conectados a um banco de dados para cadastrar os usuários." @(
    "Construir um ",
    "visual de um ",
    "site funcional de filmes, séries e animes conectados a um banco de dados para cadastrar os usuários",
    ", e exibir um gráfico dinâmico com interação ao usuario",
    "."
)

# 5) "Justificativa" paragraph content rewrite (split across 6 runs)
Set-MultiRunParagraph $d "O tema foi escolhido com base no que eu gosto que preenchesse os requisitos pedidos pelos professores" @(
    "O tema foi escolhido com base",
    " e reflexão",
    " no",
    " meu passado no que eu sentia o",
    " que eu gosto que preenchesse os requisitos pedidos pelos professores",
    " para reproduzir o projeto"
)

# 6) "Escopo" section: unchanged, nothing to do.

# 7) "Premissas e Restricoes" paragraph content rewrite (split across 5 runs)
Set-MultiRunParagraph $d "Fazer o programa funcionar da melhor forma possível tempo proposto até 05 de junho" @(
    "Fazer o programa funcionar da melhor forma possível",
    ", com todas as validações e funções acertadas até a data de apresentação",
    " tempo proposto até ",
    "14",
    " de junho"
)

Write-Output "done"
